$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Problems"): add the Problems & Solutions bullet content and
# update the title to "Problems & Solutions".
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

$content10 = $s10.Shapes.Item(1).TextFrame.TextRange

# Paragraph 1 (level 0): "SOAP for communication"
$content10.InsertAfter("SOAP ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("for")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("communication")

# Paragraph 2 (level 1): "Solution: Use RESTful service"
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("`rSolution: ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("Use")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("RESTful")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("service")
$tf10 = $s10.Shapes.Item(1).TextFrame.TextRange
$tf10.Paragraphs(2).IndentLevel = 2

# Paragraph 3 (level 0): "Update message list properly"
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("`rUpdate ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("message")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("list")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("properly")

# Paragraph 4 (level 1): "Solution: Fire DataChangeEvent"
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("`rSolution: ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("Fire")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s10.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("DataChangeEvent")
$tf10b = $s10.Shapes.Item(1).TextFrame.TextRange
$tf10b.Paragraphs(4).IndentLevel = 2

# Title: "Problems" -> "Problems & Solutions"
$s10.Shapes.Item(5).TextFrame.TextRange.Text = "Problems & Solutions"

# ---------------------------------------------------------------------------
# Slide 8 ("Server"): add the technology bullet content.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$content8 = $s8.Shapes.Item(1).TextFrame.TextRange

# Paragraph 1: "Tomcat 6"
$content8.InsertAfter("Tomcat")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" 6")

# Paragraph 2: "Jersey JAX-RS Framework"
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("`rJersey JAX-RS Framework")

# Paragraph 3: "Google Guava for Caching (with expiration)"
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("`rGoogle ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("Guava")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("for")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" Caching (")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("with")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("expiration")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(")")

# Paragraph 4: "Haversine algorithm for distance calculation"
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("`rHaversine")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("algorithm")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("for")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("distance")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter(" ")
$r = $s8.Shapes.Item(1).TextFrame.TextRange
$r.InsertAfter("calculation")

# ---------------------------------------------------------------------------
# Slide 2 ("Agenda"): merge the separate "Problems" / "Solutions" bullet
# paragraphs into a single "Problems & Solutions" paragraph.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Paragraphs(8).Delete()
$tr2b = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2b.Paragraphs(7).Characters(1, 8).Text = "Problems & Solutions"

# ---------------------------------------------------------------------------
# Slide 11 ("Solutions"): remove this now-redundant, empty slide since its
# content has been folded into slide 10's "Problems & Solutions".
# ---------------------------------------------------------------------------
$p.Slides.Item(11).Delete()
